$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 26.99883753119568
$ws.Range("C2").Value = 8.855210140476471
$ws.Range("D2").Value = 4.950014309062876
$ws.Range("E2").Value = 9.36672272388174
$ws.Range("F2").Value = 68.30292685251175
$ws.Range("J2").Value = 10.29205657185642
$ws.Range("L2").Value = 11.32679254332443
$ws.Range("M2").Value = 21.31983020993014
$ws.Range("B3").Value = 26.93403177437474
$ws.Range("C3").Value = 8.667509640607413
$ws.Range("D3").Value = 4.805650973200032
$ws.Range("E3").Value = 9.353242005042658
$ws.Range("F3").Value = 67.57095694723468
$ws.Range("J3").Value = 10.28563853136753
$ws.Range("L3").Value = 11.36752994719738
$ws.Range("M3").Value = 21.36956271853749
$ws.Range("B4").Value = 26.90388259944216
$ws.Range("C4").Value = 8.555596636862379
$ws.Range("D4").Value = 4.714771103798084
$ws.Range("E4").Value = 9.344773203179127
$ws.Range("F4").Value = 67.12492011574486
$ws.Range("J4").Value = 10.28170200439342
$ws.Range("L4").Value = 11.39432731145436
$ws.Range("M4").Value = 21.40543528162786
$ws.Range("B5").Value = 26.89402894073283
$ws.Range("C5").Value = 8.51091901792311
$ws.Range("D5").Value = 4.677211216438799
$ws.Range("E5").Value = 9.341273557804922
$ws.Range("F5").Value = 66.94413012048187
$ws.Range("J5").Value = 10.28009848419244
$ws.Range("L5").Value = 11.40569693194422
$ws.Range("M5").Value = 21.42139290878839
$ws.Range("B6").Value = 26.89253982718391
$ws.Range("C6").Value = 8.503558944484322
$ws.Range("D6").Value = 4.670943836673716
$ws.Range("E6").Value = 9.340689502033173
$ws.Range("F6").Value = 66.91417232428827
$ws.Range("J6").Value = 10.27983224811904
$ws.Range("L6").Value = 11.40761201876157
$ws.Range("M6").Value = 21.42412346553703
$ws.Range("B7").Value = 26.90373985277463
$ws.Range("C7").Value = 8.55499022549264
$ws.Range("D7").Value = 4.714266634325768
$ws.Range("E7").Value = 9.344726202757967
$ws.Range("F7").Value = 67.12247782044363
$ws.Range("J7").Value = 10.2816803766767
$ws.Range("L7").Value = 11.39447882517588
$ws.Range("M7").Value = 21.40564507281403
$ws.Range("B8").Value = 26.97449698402177
$ws.Range("C8").Value = 8.789852399693379
$ws.Range("D8").Value = 4.900724548217865
$ws.Range("E8").Value = 9.362114114223271
$ws.Range("F8").Value = 68.04991154310349
$ws.Range("J8").Value = 10.28984252091194
$ws.Range("L8").Value = 11.34046889989923
$ws.Range("M8").Value = 21.33586884454492
$ws.Range("B9").Value = 27.18926813475761
$ws.Range("C9").Value = 9.272940826567737
$ws.Range("D9").Value = 5.247080012838111
$ws.Range("E9").Value = 9.394708957835144
$ws.Range("F9").Value = 69.8901845378369
$ws.Range("J9").Value = 10.30590026089227
$ws.Range("L9").Value = 11.24868034554166
$ws.Range("M9").Value = 21.24149060556588
$ws.Range("B10").Value = 27.39254815304925
$ws.Range("C10").Value = 9.63637125822431
$ws.Range("D10").Value = 5.487965326197667
$ws.Range("E10").Value = 9.417765620330137
$ws.Range("F10").Value = 71.2481106251748
$ws.Range("J10").Value = 10.31775597773338
$ws.Range("L10").Value = 11.18980813921334
$ws.Range("M10").Value = 21.19816546756408
$ws.Range("B11").Value = 27.49465878587691
$ws.Range("C11").Value = 9.802490193877226
$ws.Range("D11").Value = 5.594271378806273
$ws.Range("E11").Value = 9.428066808892677
$ws.Range("F11").Value = 71.86558032317369
$ws.Range("J11").Value = 10.32316779890953
$ws.Range("L11").Value = 11.16487592193091
$ws.Range("M11").Value = 21.18412688121868
$ws.Range("B12").Value = 27.53468636228575
$ws.Range("C12").Value = 9.865425677904792
$ws.Range("D12").Value = 5.634031928620238
$ws.Range("E12").Value = 9.43194106613282
$ws.Range("F12").Value = 72.09923139851941
$ws.Range("J12").Value = 10.32522023419785
$ws.Range("L12").Value = 11.1556999225716
$ws.Range("M12").Value = 21.17962741773692
$ws.Range("B13").Value = 27.52600563102514
$ws.Range("C13").Value = 9.85187135914593
$ws.Range("D13").Value = 5.625491171586563
$ws.Range("E13").Value = 9.431107851751626
$ws.Range("F13").Value = 72.04892020073682
$ws.Range("J13").Value = 10.32477806111901
$ws.Range("L13").Value = 11.15766434874583
$ws.Range("M13").Value = 21.18056012082165
$ws.Range("B14").Value = 27.49792476702464
$ws.Range("C14").Value = 9.80766780163197
$ws.Range("D14").Value = 5.597552559716158
$ws.Range("E14").Value = 9.428386076934377
$ws.Range("F14").Value = 71.88480713514204
$ws.Range("J14").Value = 10.32333658952548
$ws.Range("L14").Value = 11.16411569351626
$ws.Range("M14").Value = 21.18374033489761
$ws.Range("B15").Value = 27.48090080556745
$ws.Range("C15").Value = 9.780593313405776
$ws.Range("D15").Value = 5.58037417507296
$ws.Range("E15").Value = 9.426715460846156
$ws.Range("F15").Value = 71.78425687328188
$ws.Range("J15").Value = 10.32245406228907
$ws.Range("L15").Value = 11.16810185952449
$ws.Range("M15").Value = 21.18579469056924
$ws.Range("B16").Value = 27.38606680719595
$ws.Range("C16").Value = 9.625523651884436
$ws.Range("D16").Value = 5.480949946971476
$ws.Range("E16").Value = 9.417088673912771
$ws.Range("F16").Value = 71.20774084045136
$ws.Range("J16").Value = 10.31740271275294
$ws.Range("L16").Value = 11.19147467410657
$ws.Range("M16").Value = 21.1991971388075
$ws.Range("B17").Value = 27.33034042811799
$ws.Range("C17").Value = 9.530535201172455
$ws.Range("D17").Value = 5.41909965947002
$ws.Range("E17").Value = 9.411135215614816
$ws.Range("F17").Value = 70.85390506350286
$ws.Range("J17").Value = 10.31430893642217
$ws.Range("L17").Value = 11.20628629167107
$ws.Range("M17").Value = 21.20887240952101
$ws.Range("B18").Value = 27.29919691311299
$ws.Range("C18").Value = 9.475979732572858
$ws.Range("D18").Value = 5.383218151834826
$ws.Range("E18").Value = 9.40769328245508
$ws.Range("F18").Value = 70.6503718086049
$ws.Range("J18").Value = 10.31253120320322
$ws.Range("L18").Value = 11.21497961817071
$ws.Range("M18").Value = 21.21497102568563
$ws.Range("B19").Value = 27.28880907440381
$ws.Range("C19").Value = 9.45752445657614
$ws.Range("D19").Value = 5.371017374245121
$ws.Range("E19").Value = 9.406524858460626
$ws.Range("F19").Value = 70.58146025939027
$ws.Range("J19").Value = 10.31192957672926
$ws.Range("L19").Value = 11.21795294534452
$ws.Range("M19").Value = 21.21712752355839
$ws.Range("B20").Value = 27.33617871176617
$ws.Range("C20").Value = 9.54063929237649
$ws.Range("D20").Value = 5.42571567213218
$ws.Range("E20").Value = 9.411770796634006
$ws.Range("F20").Value = 70.89157414663626
$ws.Range("J20").Value = 10.31463809259728
$ws.Range("L20").Value = 11.20469155791028
$ws.Range("M20").Value = 21.20778721868611
$ws.Range("B21").Value = 27.5061360892674
$ws.Range("C21").Value = 9.820651295967245
$ws.Range("D21").Value = 5.605772430462021
$ws.Range("E21").Value = 9.429186246581644
$ws.Range("F21").Value = 71.93301683302113
$ws.Range("J21").Value = 10.32375989728367
$ws.Range("L21").Value = 11.16221358198017
$ws.Range("M21").Value = 21.1827840571037
$ws.Range("B22").Value = 27.62512947706044
$ws.Range("C22").Value = 10.00379357003597
$ws.Range("D22").Value = 5.720553678965331
$ws.Range("E22").Value = 9.440413497980035
$ws.Range("F22").Value = 72.61260016147503
$ws.Range("J22").Value = 10.32973976651222
$ws.Range("L22").Value = 11.1359978157558
$ws.Range("M22").Value = 21.1712034055803
$ws.Range("B23").Value = 27.56090516870082
$ws.Range("C23").Value = 9.906061046902119
$ws.Range("D23").Value = 5.659565205970705
$ws.Range("E23").Value = 9.434435340030136
$ws.Range("F23").Value = 72.25003469091578
$ws.Range("J23").Value = 10.32654638898512
$ws.Range("L23").Value = 11.14984839454339
$ws.Range("M23").Value = 21.17694831830824
$ws.Range("B24").Value = 27.33353643622193
$ws.Range("C24").Value = 9.536071056064856
$ws.Range("D24").Value = 5.422725576852025
$ws.Range("E24").Value = 9.411483510320135
$ws.Range("F24").Value = 70.87454428145489
$ws.Range("J24").Value = 10.31448927825708
$ws.Range("L24").Value = 11.20541198278502
$ws.Range("M24").Value = 21.2082761635007
$ws.Range("B25").Value = 27.12310502110845
$ws.Range("C25").Value = 9.14038794276401
$ws.Range("D25").Value = 5.155645772767641
$ws.Range("E25").Value = 9.38604945364566
$ws.Range("F25").Value = 69.39077857013322
$ws.Range("J25").Value = 10.30154792998808
$ws.Range("L25").Value = 11.27200423760331
$ws.Range("M25").Value = 21.26246242929769
